# Updates the "cryptos" price table (GitHub Actions refresh) to the
# latest scraped Price / Volume(1h) figures, and fixes the row order for
# two coin pairs whose ranking swapped (NEARProtocol/Filecoin and
# Maker/Kaspa) since the previous run.
#
# Note: Price values in column D are prefixed with a leading apostrophe
# so Excel stores them as literal text instead of re-interpreting them
# as numbers (which would mangle values like "64.759.87", drop
# significant trailing zeros such as "2.80", or convert tiny decimals
# like "0.0000225" into scientific notation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.759.87"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "'3.157.82"
$ws.Range("E3").Value = "  +3.80%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'563.03"
$ws.Range("E5").Value = "  +2.52%  "
$ws.Range("D6").Value = "'146.39"
$ws.Range("E6").Value = "  +5.83%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'3.148.62"
$ws.Range("E8").Value = "  +3.82%  "
$ws.Range("D9").Value = "'0.499"
$ws.Range("E9").Value = "  +2.77%  "
$ws.Range("D10").Value = "'6.79"
$ws.Range("E10").Value = "  +5.89%  "
$ws.Range("D11").Value = "'0.156"
$ws.Range("E11").Value = "  +2.72%  "
$ws.Range("D12").Value = "'0.473"
$ws.Range("E12").Value = "  +3.51%  "
$ws.Range("D13").Value = "'37.34"
$ws.Range("E13").Value = "  +4.41%  "
$ws.Range("D14").Value = "'0.0000225"
$ws.Range("E14").Value = "  +3.13%  "
$ws.Range("D15").Value = "'3.675.98"
$ws.Range("E15").Value = "  +3.99%  "
$ws.Range("D16").Value = "'64.976.58"
$ws.Range("E16").Value = "  +1.89%  "
$ws.Range("D17").Value = "'3.177.35"
$ws.Range("E17").Value = "  +4.61%  "
$ws.Range("E18").Value = "  +2.17%  "
$ws.Range("D19").Value = "'521.48"
$ws.Range("E19").Value = "  +7.84%  "
$ws.Range("D20").Value = "'6.92"
$ws.Range("E20").Value = "  +5.51%  "
$ws.Range("D21").Value = "'14.18"
$ws.Range("E21").Value = "  +3.88%  "
$ws.Range("D22").Value = "'0.723"
$ws.Range("E22").Value = "  +5.91%  "
$ws.Range("D23").Value = "'7.54"
$ws.Range("E23").Value = "  +6.08%  "
$ws.Range("D24").Value = "'13.02"
$ws.Range("E24").Value = "  +4.77%  "
$ws.Range("D25").Value = "'79.24"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").Value = "'0.994"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").Value = "'9.09"
$ws.Range("E27").Value = "  +18.24%  "
$ws.Range("D28").Value = "'2.86"
$ws.Range("E28").Value = "  +6.07%  "
$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = "  +5.83%  "
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").Value = "'26.73"
$ws.Range("E31").Value = "  +4.09%  "
$ws.Range("D32").Value = "'2.63"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("E33").Value = "  +4.05%  "
$ws.Range("D34").Value = "'547.78"
$ws.Range("E34").Value = "  -3.84%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'5.45"
$ws.Range("E35").Value = "  +1.63%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "'6.14"
$ws.Range("E36").Value = "  +5.39%  "
$ws.Range("D37").Value = "'54.17"
$ws.Range("E37").Value = "  +4.80%  "
$ws.Range("D38").Value = "'0.0438"
$ws.Range("E38").Value = "  +6.80%  "
$ws.Range("D39").Value = "'0.0833"
$ws.Range("E39").Value = "  +5.20%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "'3.169.45"
$ws.Range("E40").Value = "  +8.59%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.124"
$ws.Range("E41").Value = "  +5.64%  "
$ws.Range("D42").Value = "'2.80"
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("D43").Value = "'8.37"
$ws.Range("E43").Value = "  +2.07%  "
$ws.Range("D44").Value = "'0.272"
$ws.Range("E44").Value = "  +12.59%  "
$ws.Range("D45").Value = "'2.24"
$ws.Range("E45").Value = "  +8.61%  "
$ws.Range("D47").Value = "'25.84"
$ws.Range("E47").Value = "  +4.42%  "
$ws.Range("D48").Value = "'122.95"
$ws.Range("E48").Value = "  +4.30%  "
$ws.Range("D49").Value = "'0.0₃0530"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("E50").Value = "  +1.35%  "
$ws.Range("D51").Value = "'2.13"
$ws.Range("E51").Value = "  +4.48%  "
